$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Report")

$updates = @{
    "AW2" = 157.855949
    "AK3" = 48.694259
    "AQ4" = 20.645301
    "AW5" = 0.971944
    "AW6" = 112.801609
    "AW7" = 115.58816
    "AK8" = 15.693275
    "AW9" = 150.62566
    "AW10" = 99.666296
    "AW11" = 83.749583
    "AW12" = 84.651956
    "AQ13" = 0.974479
    "AW14" = 162.709653
    "AW15" = 148.524803
    "AW16" = 117.653333
    "AW17" = 119.950637
    "AN18" = 22.823669
    "AQ19" = 10.387106
    "AW20" = 87.537465
    "AW21" = 79.626227
    "AK22" = 15.693287
    "AK23" = 15.693333
    "AW24" = 0.972373
    "AW25" = 112.801076
    "AQ26" = 1.645185
    "AK27" = 15.693438
    "AW28" = 148.530313
    "AW29" = 69.613796
    "AW30" = 78.70728
    "AW31" = 8.847326
    "AW32" = 150.620567
    "AW33" = 86.849282
    "AW34" = 7.628727
    "AQ35" = 1.645926
    "AW36" = 114.895081
    "AW37" = 122.680822
    "AW38" = 9.91272
    "AK39" = 15.69331
    "AK40" = 15.693206
    "AK41" = 15.693519
    "AK42" = 15.693229
    "AQ43" = 20.645174
    "AW44" = 71.932083
    "AW45" = 86.848669
    "AW46" = 37.96544
    "AW47" = 37.964225
    "AK48" = 37.742824
    "AW49" = 114.896065
    "AW50" = 115.589572
    "AW51" = 22.732188
    "AW52" = 128.805058
    "AW53" = 106.562593
    "AW54" = 153.667245
    "AW55" = 143.716065
    "AK56" = 38.543067
    "AW57" = 115.588738
    "AQ58" = 20.644688
    "AW59" = 9.916829
    "AW60" = 51.590313
    "AW61" = 118.844421
    "AW62" = 79.645567
    "AK63" = 10.719039
    "AQ64" = 2.911829
    "AW65" = 154.960382
    "AW66" = 157.853275
    "AW67" = 115.604514
    "AW68" = 72.706042
    "AW69" = 15.7386
    "AK70" = 15.693495
    "AK71" = 15.693252
    "AW72" = 8.742639
    "AK73" = 37.742905
    "AW74" = 129.675949
    "AW75" = 66.702164
    "AW76" = 1.629155
    "AW77" = 70.91081
    "AN78" = 8.858808
    "AW79" = 143.719456
    "AW80" = 37.96522
    "AW81" = 22.880579
    "AW82" = 17.575336
    "AW83" = 7.627384
    "AW84" = 122.74478
    "AW85" = 153.619525
    "AW86" = 101.58919
    "AW87" = 70.909097
    "AW88" = 51.58316
    "AQ89" = 35.674965
    "AW90" = 107.591701
    "AW91" = 148.524421
    "AW92" = 57.96485
    "AK93" = 10.716597
    "AW94" = 115.587963
    "AK95" = 15.69316
    "AW96" = 154.960775
    "AW97" = 66.57544
    "AQ98" = 20.645405
    "AW99" = 84.644618
    "AW100" = 76.720475
    "AW101" = 150.710382
    "AW102" = 106.707373
    "AW103" = 45.966123
    "AK104" = 48.694225
    "AK105" = 10.719063
    "AQ106" = 3.871204
    "AW107" = 105.144248
    "AW108" = 113.787766
    "AW109" = 66.631389
    "AK110" = 15.69338
    "AN111" = 8.858762
    "AW112" = 15.720081
    "AK113" = 15.693194
    "AW114" = 157.856238
    "AW115" = 143.667118
    "AW116" = 153.617662
    "AW117" = 45.964525
    "AW118" = 17.565972
    "AQ119" = 1.647766
    "AW120" = 132.586331
    "AW121" = 106.727384
    "AW122" = 84.597824
    "AW123" = 70.906944
    "AK124" = 15.693322
    "AK125" = 38.543171
    "AW126" = 148.525035
    "AK127" = 15.693229
    "AN128" = 8.858831
    "AW129" = 143.61963
    "AK130" = 37.74294
    "AK131" = 15.693206
    "AW132" = 57.827199
    "AW133" = 51.978611
    "AW134" = 161.669653
    "AW135" = 87.510012
    "AW136" = 45.963981
    "AK137" = 15.693299
    "AN138" = 22.748542
    "AW139" = 143.670324
    "AW140" = 148.52669
    "AW141" = 84.707315
    "AK142" = 15.693241
    "AW143" = 118.843993
    "AW144" = 115.588461
    "AW145" = 78.70831
    "AK146" = 15.693484
    "AW147" = 50.840486
    "AW148" = 50.640544
    "AQ149" = 35.67485
    "AW150" = 72.538681
    "AW151" = 62.603912
    "AW152" = 92.674086
    "AW153" = 45.963287
    "AK154" = 15.693345
    "AQ155" = 27.654965
    "AW156" = 114.895868
    "AW157" = 70.912558
    "AK158" = 10.716597
    "AN159" = 3.726065
    "AW160" = 104.918681
    "AW161" = 20.665127
    "AN162" = 0.54544
    "AW163" = 153.620463
    "AW164" = 129.970764
    "AW165" = 69.613588
    "AK166" = 15.693171
    "AW167" = 104.923623
    "AW168" = 157.856632
    "AW169" = 79.624699
    "AW170" = 150.620012
    "AW171" = 141.726852
    "AW172" = 79.645451
    "AW173" = 8.742477
    "AK174" = 37.742813
    "AW175" = 150.619653
    "AW176" = 106.723194
    "AW177" = 51.977894
    "AW178" = 84.537604
    "AN179" = 7.690544
    "AW180" = 157.864201
    "AW181" = 92.686273
    "AW182" = 70.910069
    "AW183" = 51.981435
    "AW184" = 106.727014
    "AK185" = 15.693322
    "AK186" = 37.74294
    "AW187" = 150.71066
    "AW188" = 164.931493
    "AW189" = 70.915914
    "AW190" = 44.70331
    "AW191" = 137.947199
    "AW192" = 148.525602
    "AW193" = 76.629155
    "AW194" = 132.583171
    "AW195" = 80.552593
    "AK196" = 10.719051
    "AW197" = 153.593796
    "AW198" = 150.710116
    "AW199" = 149.632338
    "AW200" = 161.614919
    "AW201" = 163.755313
    "AW202" = 153.668495
    "AW203" = 115.590185
    "AW204" = 20.933194
    "AK205" = 37.742824
    "AW206" = 115.60816
    "AW207" = 162.710637
    "AW208" = 153.60956
    "AW209" = 162.706505
    "AW210" = 63.66081
    "AK211" = 15.693275
    "AK212" = 15.693264
    "AQ213" = 20.644907
    "AN214" = 7.690521
    "AW215" = 70.915266
    "AW216" = 122.680741
    "AW217" = 139.922049
    "AW218" = 79.691655
    "AW219" = 71.931852
    "AW220" = 51.979745
    "AW221" = 2.908611
    "AW222" = 114.8964
    "AW223" = 70.774861
    "AK224" = 37.742928
    "AQ225" = 27.647708
    "AW226" = 104.922396
    "AW227" = 150.594896
    "AW228" = 115.58897
    "AW229" = 79.623044
    "AW230" = 22.877627
    "AW231" = 51.682512
    "AW232" = 45.962407
    "AK233" = 15.693218
    "AK234" = 38.543183
    "AW235" = 115.589375
    "AW236" = 86.849907
    "AK237" = 15.693299
    "AW238" = 135.742569
    "AK239" = 15.686296
    "AK240" = 15.693171
    "AW241" = 114.856979
    "AW242" = 45.972153
    "AW243" = 45.97
    "AW244" = 20.554549
    "AK245" = 15.693252
    "AW246" = 112.801238
    "AW247" = 36.960718
    "AW248" = 22.73294
    "AW249" = 143.719676
    "AW250" = 153.667755
    "AW251" = 132.583519
    "AW252" = 22.73375
    "AN253" = 8.858785
    "AW254" = 8.742373
    "AW255" = 132.584063
    "AW256" = 51.974144
    "AW257" = 3.720984
    "AW258" = 142.61287
    "AW259" = 57.828391
    "AW260" = 9.915
    "AW261" = 128.804861
    "AQ262" = 20.644815
    "AW263" = 8.848056
    "AQ264" = 2.907755
    "AW265" = 37.964144
    "AQ266" = 2.907639
    "AW267" = 129.971134
    "AW268" = 141.727801
    "AW269" = 104.923461
    "AW270" = 84.651528
    "AW271" = 25.978947
    "AW272" = 128.93353
    "AK273" = 37.742836
    "AW274" = 3.714641
    "AW275" = 142.932731
    "AW276" = 87.540625
    "AW277" = 148.523634
    "AW278" = 70.916481
    "AW279" = 115.591308
    "AK280" = 66.75912
    "AW281" = 139.922512
    "AW282" = 70.909352
    "AW283" = 45.720822
    "AK284" = 15.693426
    "AW285" = 80.682292
    "AW286" = 154.962234
    "AW287" = 44.634444
    "AK288" = 15.693495
    "AK289" = 15.693183
    "AQ290" = 1.644838
    "AW291" = 44.633009
    "AW292" = 44.63478
    "AK293" = 38.543183
    "AW294" = 148.524016
    "AW295" = 101.585706
    "AW296" = 20.666076
    "AW297" = 101.588495
    "AW298" = 114.893738
    "AK299" = 10.716586
    "AW300" = 157.857083
    "AW301" = 51.973287
    "AW302" = 15.683368
    "AW303" = 3.718623
    "AW304" = 136.631134
    "AW305" = 91.69015
    "AK306" = 15.693403
    "AK307" = 15.693356
    "AW308" = 115.590637
    "AW309" = 84.597523
    "AW310" = 30.608241
    "AW311" = 9.911215
    "AQ312" = 2.911898
    "AW313" = 94.615116
    "AK314" = 15.693426
    "AW315" = 86.846551
    "AW316" = 20.66559
    "AW317" = 162.715463
    "AW318" = 72.914433
    "AW319" = 51.980648
    "AW320" = 1.628669
    "AW321" = 119.953808
    "AK322" = 15.693356
    "AK323" = 15.693403
    "AQ324" = 2.907894
    "AW325" = 76.67419
    "AW326" = 55.612708
    "AW327" = 20.666701
    "AW328" = 1.629896
    "AW329" = 128.93375
    "AW330" = 57.82794
    "AW331" = 45.961644
    "AW332" = 52.034155
    "AK333" = 15.693472
    "AQ334" = 3.871308
    "AW335" = 153.618368
    "AW336" = 154.960544
    "AK337" = 15.693368
    "AK338" = 23.644803
    "AW339" = 161.609213
    "AW340" = 51.981863
    "AW341" = 132.584583
    "AW342" = 22.879132
    "AW343" = 20.914167
    "AW344" = 148.523368
    "AW345" = 119.951123
    "AW346" = 38.551319
    "AW347" = 122.743681
    "AW348" = 70.910706
    "AW349" = 157.74912
    "AW350" = 153.654977
    "AW351" = 1.635289
    "AW352" = 154.961343
    "AW353" = 132.585104
    "AW354" = 162.706863
    "AW355" = 10.643241
    "AW356" = 148.521898
    "AW357" = 115.587778
    "AW358" = 106.720405
    "AW359" = 22.878449
    "AK360" = 15.693391
    "AK361" = 37.742917
    "AW362" = 69.613935
    "AW363" = 69.613472
    "AK364" = 15.693507
    "AQ365" = 20.645058
    "AW366" = 1.628877
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

$wb.Save()